{"js": "// The prtGen installation instructions switched from shell-script launchers\n// (\"setup.command\" / \"run.command\") to Python launchers (\"setup.py\" /\n// \"prtGen.py\"). Update the two distinct pieces of text throughout the body.\n\nconst body = context.document.body;\n\n// 1) \"Double click on setup.command to install the required components\"\n//    -> \"Double click on setup.py to install the required components\"\nconst setupHits = body.search(\"setup.command\", { matchCase: true });\nsetupHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < setupHits.items.length; i++) {\n  setupHits.items[i].insertText(\"setup.py\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"run.command\" -> \"prtGen.py\" (appears twice: once describing how to\n//    start the application, once describing how to create a shortcut/alias).\nconst runHits = body.search(\"run.command\", { matchCase: true });\nrunHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < runHits.items.length; i++) {\n  runHits.items[i].insertText(\"prtGen.py\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The prtGen installation instructions switched from shell-script launchers\n# (\"setup.command\" / \"run.command\") to Python launchers (\"setup.py\" /\n# \"prtGen.py\"). Update both pieces of text throughout the document body.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) \"Double click on setup.command to install the required components\"\n#    -> \"Double click on setup.py to install the required components\"\nReplace-AllText \"setup.command\" \"setup.py\"\n\n# 2) \"run.command\" -> \"prtGen.py\" (appears twice: once describing how to\n#    start the application, once describing how to create a shortcut/alias).\nReplace-AllText \"run.command\" \"prtGen.py\"\n"}
